$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = -0.3030476147762293
$ws.Range("C2").Value = -0.2852167463729813
$ws.Range("B3").Value = 0.2212987470419176
$ws.Range("C3").Value = 0.2314516030111665
$ws.Range("B4").Value = 0.1635776920988358
$ws.Range("C4").Value = 0.1762598686712164
$ws.Range("B5").Value = -0.04981953747724374
$ws.Range("C5").Value = -0.04037048095642987
$ws.Range("B6").Value = 0.1575450856876263
$ws.Range("C6").Value = 0.1734455518791605
$ws.Range("B7").Value = -0.4379502147954654
$ws.Range("C7").Value = -0.4275206503871763
$ws.Range("B8").Value = -0.2808802019124831
$ws.Range("C8").Value = -0.2590080310344377
$ws.Range("B9").Value = -0.4047525790627515
$ws.Range("C9").Value = -0.3942749422504559
$ws.Range("B10").Value = 0.3294102058557918
$ws.Range("C10").Value = 0.3468392264229247
$ws.Range("B11").Value = -0.202385559409978
$ws.Range("C11").Value = -0.1944054799455128
$ws.Range("B12").Value = -0.08263167799945226
$ws.Range("C12").Value = -0.05724445409281485
$ws.Range("B13").Value = 0.003378030045047247
$ws.Range("C13").Value = 0.003991582187646852
$ws.Range("B14").Value = 0.03887108225803239
$ws.Range("C14").Value = 0.07897491892330588
$ws.Range("B15").Value = -0.07872759690180693
$ws.Range("C15").Value = -0.02953390677198398
$ws.Range("B16").Value = 0.2669074388759973
$ws.Range("C16").Value = 0.3329519714974425
$ws.Range("B17").Value = 0.5781798699301779
$ws.Range("C17").Value = 0.6104878113432393
$ws.Range("B18").Value = 0.05681595236303819
$ws.Range("C18").Value = 0.03871942535414385
$ws.Range("B19").Value = 0.4004881803427321
$ws.Range("C19").Value = 0.4094959229080015
$ws.Range("B20").Value = 0.2416489704046453
$ws.Range("C20").Value = 0.3007323468731548
$ws.Range("B21").Value = 0.4429884574041729
$ws.Range("C21").Value = 0.5097515718032218
$ws.Range("B22").Value = 0.3506008479462824
$ws.Range("C22").Value = 0.3869851860350284
$ws.Range("B23").Value = -0.06633768584288893
$ws.Range("C23").Value = -0.03439731083199027
$ws.Range("B24").Value = 4.539342217017882
$ws.Range("C24").Value = 4.562109517439753
$ws.Range("B25").Value = 0.542396028553592
$ws.Range("C25").Value = 0.5076937306151328
$ws.Range("B26").Value = 0.4019698531004052
$ws.Range("C26").Value = 0.3829019154812469
$ws.Range("B27").Value = 0.3324666689593566
$ws.Range("C27").Value = 0.3026541687337232
$ws.Range("B28").Value = 1.102688686605339
$ws.Range("C28").Value = 1.07192136330336
$ws.Range("B29").Value = 5.829782112900279
$ws.Range("C29").Value = 5.37474763382918
$ws.Range("B30").Value = 1.002976036868044
$ws.Range("C30").Value = 0.9535325995048894
$ws.Range("B31").Value = -0.1433614335482137
$ws.Range("C31").Value = -0.2065261134623893
$ws.Range("B32").Value = 0.8203836391126612
$ws.Range("C32").Value = 0.7801732546091538
$ws.Range("B33").Value = 0.9197923190188908
$ws.Range("C33").Value = 0.89277604747635
$ws.Range("B34").Value = -0.5953512697313745
$ws.Range("C34").Value = -0.6266926107944895
$ws.Range("B35").Value = 0.8319897273086079
$ws.Range("C35").Value = 0.8220127442528157
$ws.Range("B36").Value = 0.7852938092493125
$ws.Range("C36").Value = 0.7691359516405509
$ws.Range("B37").Value = 0.767646983259594
$ws.Range("C37").Value = 0.7471079223737263
$ws.Range("B38").Value = 0.759471446328237
$ws.Range("C38").Value = 0.7365332892878046
$ws.Range("B39").Value = 0.5831043790969885
$ws.Range("C39").Value = 0.5807572289649692
$ws.Range("B40").Value = 0.7545904765981907
$ws.Range("C40").Value = 0.7531474682819043
$ws.Range("B41").Value = 0.5736631958143509
$ws.Range("C41").Value = 0.5660166821978613
$ws.Range("B42").Value = 0.7215599402919902
$ws.Range("C42").Value = 0.6920796923304855
$ws.Range("B43").Value = 0.7336091938271017
$ws.Range("C43").Value = 0.7179243178824716
$ws.Range("B44").Value = 0.6743193507660545
$ws.Range("C44").Value = 0.6672284289971849
$ws.Range("B45").Value = 0.6600740514179404
$ws.Range("C45").Value = 0.6429950345313618
$ws.Range("B46").Value = -1.252431146668739
$ws.Range("C46").Value = -1.256512184098387
$ws.Range("B47").Value = -0.9698305569425633
$ws.Range("C47").Value = -0.9747869049202756
$ws.Range("B48").Value = -0.8627063238443551
$ws.Range("C48").Value = -0.8687641519963322
$ws.Range("B49").Value = -0.6316061218671422
$ws.Range("C49").Value = -0.6349765219189132
$ws.Range("B50").Value = -0.04877058266929571
$ws.Range("C50").Value = -0.049114311899176
$ws.Range("B51").Value = -0.8513032994962222
$ws.Range("C51").Value = -0.8551400172311395
$ws.Range("B52").Value = -0.8513032994962222
$ws.Range("C52").Value = -0.8551400172311395
$ws.Range("B53").Value = -1.079134497588794
$ws.Range("C53").Value = -1.093551729333249
$ws.Range("B54").Value = -0.1865897922561882
$ws.Range("C54").Value = -0.184982862779744
$ws.Range("B55").Value = -0.9900142941882175
$ws.Range("C55").Value = -0.9935862167009119
$ws.Range("B56").Value = -0.8920081595603776
$ws.Range("C56").Value = -0.8860932383385447
$ws.Range("B57").Value = -0.9651520261654113
$ws.Range("C57").Value = -0.9480444094739415
$ws.Range("B58").Value = -1.165967125806632
$ws.Range("C58").Value = -1.135323957330342
$ws.Range("B59").Value = -0.8696744983449465
$ws.Range("C59").Value = -0.8518372267606601
$ws.Range("B60").Value = -0.5216568842680853
$ws.Range("C60").Value = -0.4958122221176305
$ws.Range("B61").Value = 0.3668232733977567
$ws.Range("C61").Value = 0.3699346703233274
$ws.Range("B62").Value = -1.255433272480964
$ws.Range("C62").Value = -1.232361134040926
$ws.Range("B63").Value = -0.762466120581034
$ws.Range("C63").Value = -0.7268749157402274
$ws.Range("B64").Value = -0.9010610062281696
$ws.Range("C64").Value = -0.8924581126356002
$ws.Range("B65").Value = -0.1406406803056653
$ws.Range("C65").Value = -0.115100628572694
$ws.Range("B66").Value = -0.8330872961444159
$ws.Range("C66").Value = -0.8023975966178745
$ws.Range("B67").Value = -0.8361999719167592
$ws.Range("C67").Value = -0.7939409826044265
